# Updated TPM-derived NATMI metrics for the Ifnb1-Ifnar2 LR-pair sheet.
# Only the numeric columns E-J and M-T (expression/specificity metrics) change;
# Sending/Ligand/Receptor/Target cluster text columns (A-D) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03958166666666667
$ws.Range("H2").Value = 0.118745
$ws.Range("I2").Value = 0.2870281964201545
$ws.Range("J2").Value = 0.2870281964201545
$ws.Range("M2").Value = 31.743396
$ws.Range("N2").Value = 95.230188
$ws.Range("O2").Value = 0.1189237443612096
$ws.Range("P2").Value = 0.1189237443612096
$ws.Range("Q2").Value = 1.25645651934
$ws.Range("R2").Value = 11.30810867406
$ws.Range("S2").Value = 0.0341344678555295
$ws.Range("T2").Value = 0.0341344678555295
# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03958166666666667
$ws.Range("H3").Value = 0.118745
$ws.Range("I3").Value = 0.2870281964201545
$ws.Range("J3").Value = 0.2870281964201545
$ws.Range("O3").Value = 0.1062760886263749
$ws.Range("P3").Value = 0.106276088626375
$ws.Range("Q3").Value = 1.122831147991667
$ws.Range("R3").Value = 10.105480331925
$ws.Range("S3").Value = 0.03050423404101689
$ws.Range("T3").Value = 0.0305042340410169
# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03958166666666667
$ws.Range("H4").Value = 0.118745
$ws.Range("I4").Value = 0.2870281964201545
$ws.Range("J4").Value = 0.2870281964201545
$ws.Range("M4").Value = 70.51016133333333
$ws.Range("N4").Value = 211.530484
$ws.Range("O4").Value = 0.264159902780187
$ws.Range("P4").Value = 0.264159902780187
$ws.Range("Q4").Value = 2.790909702508889
$ws.Range("R4").Value = 25.11818732258
$ws.Range("S4").Value = 0.07582134046152042
$ws.Range("T4").Value = 0.07582134046152043
# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03958166666666667
$ws.Range("H5").Value = 0.118745
$ws.Range("I5").Value = 0.2870281964201545
$ws.Range("J5").Value = 0.2870281964201545
$ws.Range("M5").Value = 6.595865666666666
$ws.Range("N5").Value = 19.787597
$ws.Range("O5").Value = 0.02471081047483217
$ws.Range("P5").Value = 0.02471081047483218
$ws.Range("Q5").Value = 0.2610753561961111
$ws.Range("R5").Value = 2.349678205765
$ws.Range("S5").Value = 0.007092699362671339
$ws.Range("T5").Value = 0.007092699362671341
# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.03958166666666667
$ws.Range("H6").Value = 0.118745
$ws.Range("I6").Value = 0.2870281964201545
$ws.Range("J6").Value = 0.2870281964201545
$ws.Range("M6").Value = 64.15060166666667
$ws.Range("N6").Value = 192.451805
$ws.Range("O6").Value = 0.2403343912297365
$ws.Range("P6").Value = 0.2403343912297365
$ws.Range("Q6").Value = 2.539187731636111
$ws.Range("R6").Value = 22.852689584725
$ws.Range("S6").Value = 0.06898274685240705
$ws.Range("T6").Value = 0.06898274685240706
# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.03958166666666667
$ws.Range("H7").Value = 0.118745
$ws.Range("I7").Value = 0.2870281964201545
$ws.Range("J7").Value = 0.2870281964201545
$ws.Range("M7").Value = 65.55479199999999
$ws.Range("N7").Value = 196.664376
$ws.Range("O7").Value = 0.2455950625276598
$ws.Range("P7").Value = 0.2455950625276598
$ws.Range("Q7").Value = 2.594767925346666
$ws.Range("R7").Value = 23.35291132812
$ws.Range("S7").Value = 0.07049270784700924
$ws.Range("T7").Value = 0.07049270784700926
# Row 8
$ws.Range("G8").Value = 0.09832
$ws.Range("H8").Value = 0.29496
$ws.Range("I8").Value = 0.7129718035798456
$ws.Range("J8").Value = 0.7129718035798456
$ws.Range("M8").Value = 31.743396
$ws.Range("N8").Value = 95.230188
$ws.Range("O8").Value = 0.1189237443612096
$ws.Range("P8").Value = 0.1189237443612096
$ws.Range("Q8").Value = 3.12101069472
$ws.Range("R8").Value = 28.08909625248
$ws.Range("S8").Value = 0.08478927650568006
$ws.Range("T8").Value = 0.08478927650568008
# Row 9
$ws.Range("G9").Value = 0.09832
$ws.Range("H9").Value = 0.29496
$ws.Range("I9").Value = 0.7129718035798456
$ws.Range("J9").Value = 0.7129718035798456
$ws.Range("O9").Value = 0.1062760886263749
$ws.Range("P9").Value = 0.106276088626375
$ws.Range("Q9").Value = 2.7890881756
$ws.Range("R9").Value = 25.1017935804
$ws.Range("S9").Value = 0.07577185458535804
$ws.Range("T9").Value = 0.07577185458535807
# Row 10
$ws.Range("G10").Value = 0.09832
$ws.Range("H10").Value = 0.29496
$ws.Range("I10").Value = 0.7129718035798456
$ws.Range("J10").Value = 0.7129718035798456
$ws.Range("M10").Value = 70.51016133333333
$ws.Range("N10").Value = 211.530484
$ws.Range("O10").Value = 0.264159902780187
$ws.Range("P10").Value = 0.264159902780187
$ws.Range("Q10").Value = 6.932559062293334
$ws.Range("R10").Value = 62.39303156064
$ws.Range("S10").Value = 0.1883385623186666
$ws.Range("T10").Value = 0.1883385623186666
# Row 11
$ws.Range("G11").Value = 0.09832
$ws.Range("H11").Value = 0.29496
$ws.Range("I11").Value = 0.7129718035798456
$ws.Range("J11").Value = 0.7129718035798456
$ws.Range("M11").Value = 6.595865666666666
$ws.Range("N11").Value = 19.787597
$ws.Range("O11").Value = 0.02471081047483217
$ws.Range("P11").Value = 0.02471081047483218
$ws.Range("Q11").Value = 0.6485055123466666
$ws.Range("R11").Value = 5.83654961112
$ws.Range("S11").Value = 0.01761811111216083
$ws.Range("T11").Value = 0.01761811111216084
# Row 12
$ws.Range("G12").Value = 0.09832
$ws.Range("H12").Value = 0.29496
$ws.Range("I12").Value = 0.7129718035798456
$ws.Range("J12").Value = 0.7129718035798456
$ws.Range("M12").Value = 64.15060166666667
$ws.Range("N12").Value = 192.451805
$ws.Range("O12").Value = 0.2403343912297365
$ws.Range("P12").Value = 0.2403343912297365
$ws.Range("Q12").Value = 6.307287155866668
$ws.Range("R12").Value = 56.7655844028
$ws.Range("S12").Value = 0.1713516443773294
$ws.Range("T12").Value = 0.1713516443773294
# Row 13
$ws.Range("G13").Value = 0.09832
$ws.Range("H13").Value = 0.29496
$ws.Range("I13").Value = 0.7129718035798456
$ws.Range("J13").Value = 0.7129718035798456
$ws.Range("M13").Value = 65.55479199999999
$ws.Range("N13").Value = 196.664376
$ws.Range("O13").Value = 0.2455950625276598
$ws.Range("P13").Value = 0.2455950625276598
$ws.Range("Q13").Value = 6.44534714944
$ws.Range("R13").Value = 58.00812434496
$ws.Range("S13").Value = 0.1751023546806505
$ws.Range("T13").Value = 0.1751023546806506
